$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The dropdown-file questions (rows 6-9 / spreadsheet rows 7-10) referenced a
# CSV source file for the native-language dropdown lists. The source file is
# actually an .xlsx workbook, so fix the filename reference in column C.
$ws.Range("C7").Value = "language_iso639_1_toy.xlsx"
$ws.Range("C8").Value = "language_iso639_1_toy.xlsx"
$ws.Range("C9").Value = "language_iso639_1_toy.xlsx"
$ws.Range("C10").Value = "language_iso639_1_toy.xlsx"

# Reflect the last on-screen scroll position / active selection from the
# authoring session (scrolled right to column H, cell N19 selected).
$ws.Application.Goto($ws.Range("N19"), $false)
$ws.Range("N19").Select()
$excel.ActiveWindow.ScrollColumn = 8
